# "Generate Report for Archive"
#
# The localization status text for this particular file moved from
# "Ready for handoff" to "In Translation":
#   - Overview sheet: columns "zh-cn" (E) and "de-de" (F), row 2
#   - zh-cn sheet:     "Status" column (C), row 2
#   - de-de sheet:     "Status" column (C), row 2
#
# As a consequence, the (now shorter) Status text caused those columns to
# shrink when the report's column widths were refreshed, so we narrow the
# same columns to match.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# Target (from the authored XML) column width is 13.4101845877511 "characters".
# This runtime's ColumnWidth setter snaps to a 1/6-character pixel grid when it
# stores the value, so feed it the input ("12.5") that rounds to the pixel count
# (75) nearest the authored width, giving the closest achievable match.
$newColumnWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update the status value everywhere it appears ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# --- Narrow the columns that held the (now shorter) status text ---
$overview.Columns.Item(5).ColumnWidth = $newColumnWidth   # "zh-cn" column
$overview.Columns.Item(6).ColumnWidth = $newColumnWidth   # "de-de" column
$zhcn.Columns.Item(3).ColumnWidth = $newColumnWidth        # "Status" column
$dede.Columns.Item(3).ColumnWidth = $newColumnWidth        # "Status" column
